$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet "survey": add `hint` and `choice_filter` columns, and insert a
# new cascading-CSV example (select_one regions_csv / select_one
# countries_csv) just before the content-provider example, which gets
# pushed down two rows and gains a `hint` plus an updated callback in
# the queries sheet.
# -------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("H1").Value = "hint"
$survey.Range("I1").Value = "choice_filter"

# Move the existing "content_provider_test" row (currently row 17) down
# to row 19, and give it the new hint text.
$survey.Range("B19").Value = $survey.Range("B17").Value2
$survey.Range("F19").Value = $survey.Range("F17").Value2
$survey.Range("G19").Value = $survey.Range("G17").Value2
$survey.Range("H19").Value = "You will need to install a content provider app for the query to work. There is an example app available here: https://github.com/nathanathan/FileContentProviderExample"

# Row 17 becomes the new "region" cascading select.
$survey.Range("B17").Value = "select_one regions_csv"
$survey.Range("C17").Value = ""
$survey.Range("F17").Value = "region"
$survey.Range("G17").Value = "Choose a region:"

# Row 18 is brand new: the "country" cascading select, filtered by region.
$survey.Range("B18").Value = "select_one countries_csv"
$survey.Range("C18").Value = "dropdown"
$survey.Range("F18").Value = "country"
$survey.Range("G18").Value = "Choose a country:"
$survey.Range("I18").Value = "_.where(context, {`n  region: data('region')`n})"

# Keep the trailing block of blank, formatted placeholder rows the same
# size as before by extending it with two more rows (24 and 25), copying
# the formatting already used on column A of the blank rows.
$survey.Range("A23").Copy() | Out-Null
$survey.Range("A24:A25").PasteSpecial(-4122) | Out-Null

# -------------------------------------------------------------------
# Sheet "queries": widen the callback column, add queries backing the
# two new CSV-based selects, and point the content-provider query at
# the new example content provider / callback.
# -------------------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")

$queries.Columns.Item(3).ColumnWidth = 42.8

$queries.Range("A4").Value = "countries_csv"
$queries.Range("B4").Value = '"regions.csv"'
$queries.Range("C4").Value = "_.chain(context).pluck('region').uniq().map(function(region){`nreturn {name:region, label:region};`n})"

$queries.Range("A5").Value = "regions_csv"
$queries.Range("B5").Value = '"regions.csv"'
$queries.Range("C5").Value = "_.map(context, function(place){`nplace.name = place.country;`nplace.label = place.country;`nreturn place;`n})"

$queries.Range("A6").Value = "content_provider_test"
$queries.Range("B6").Value = '"content://org.opendatakit.FileContentProviderExample/"'
$queries.Range("C6").Value = "[context]"
